$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: A4/D4 switch from text "12"/"21" to numeric 12/21
$ws.Range("A4").Value = 12
$ws.Range("D4").Value = 21

# Row 5
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "10:18"
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = "Packing"
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 12
$ws.Range("J5").Value = 9694.08
$ws.Range("K5").Value = "testing purposes"

# Row 6
$ws.Range("A6").Value = 12
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "10:18"
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = "Packing"
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 12
$ws.Range("J6").Value = 9694.08
$ws.Range("K6").Value = "testing purposes"

# Row 7
$ws.Range("A7").Value = 12
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "10:19"
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = "Packing"
$ws.Range("F7").Value = 12
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = 12345
$ws.Range("I7").Value = 12
$ws.Range("J7").Value = 10259502363
$ws.Range("K7").Value = "testing purposes"

# Row 8
$ws.Range("A8").Value = "A1234P"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "10:22"
$ws.Range("D8").Value = "faris"
$ws.Range("E8").Value = "Packing"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 5.61
$ws.Range("K8").Value = "testing purposes"

# Row 9
$ws.Range("A9").Value = "A12212P"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "10:58"
$ws.Range("D9").Value = "FARIS"
$ws.Range("E9").Value = "'70"
$ws.Range("E9").ClearFormats()
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 2.805
$ws.Range("K9").Value = "lakukan pemanasan hingga 240°C"

# Row 10
$ws.Range("A10").Value = "A12212P"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "10:59"
$ws.Range("D10").Value = "FARIS"
$ws.Range("E10").Value = "'240"
$ws.Range("E10").ClearFormats()
$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 2.805
$ws.Range("K10").Value = "Lakukan cooling hingga 120°C"

# Row 11
$ws.Range("A11").Value = "A12212P"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "10:59"
$ws.Range("D11").Value = "FARIS"
$ws.Range("E11").Value = "'240"
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 28.05
$ws.Range("K11").Value = "Tambah waktu pemanasan"

# Row 12
$ws.Range("A12").Value = "A12212P"
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = "10:59"
$ws.Range("D12").Value = "FARIS"
$ws.Range("E12").Value = "'120"
$ws.Range("E12").ClearFormats()
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 50
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 28.05
$ws.Range("K12").Value = "Hubungi atasan"

# Row 13
$ws.Range("A13").Value = "A12212P"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "10:59"
$ws.Range("D13").Value = "FARIS"
$ws.Range("E13").Value = "'120"
$ws.Range("E13").ClearFormats()
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 5
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 2.805
$ws.Range("K13").Value = "Tambah Oleic Acid"

# Row 14
$ws.Range("A14").Value = "'1212"
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "11:20"
$ws.Range("D14").Value = "'12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'70"
$ws.Range("E14").ClearFormats()
$ws.Range("F14").Value = 12
$ws.Range("G14").Value = 12
$ws.Range("H14").Value = 12
$ws.Range("I14").Value = 12
$ws.Range("J14").Value = 9694.08
$ws.Range("K14").Value = "lakukan pemanasan hingga 240°C"

# Row 15
$ws.Range("A15").Value = "'12"
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "11:30"
$ws.Range("D15").Value = "'12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'40"
$ws.Range("E15").ClearFormats()
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 12
$ws.Range("I15").Value = 12
$ws.Range("J15").Value = 9694.08
$ws.Range("K15").Value = "lakukan pemanasan hingga 240°C"

# Row 16
$ws.Range("A16").Value = "'12"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = "11:30"
$ws.Range("D16").Value = "'12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "Packing"
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 12
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 9694.08
$ws.Range("K16").Value = "NG"
